# Passive Skills and Passive Skills Logic
# Fix three typos in the StatusDescription column (C) of the StatusEffect
# Database sheet, and restore the editor's last selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 (Barrier): "na atack" -> "an atack"
$ws.Range("C27").Value = "Consume a stack to block an atack"

# Row 40 (Overweigth): "dels" -> "deals"
$ws.Range("C40").Value = "Can not get extra turn, deals minus 50% damage and take plus 50% damage"

# Row 43 (Healing): "teh" -> "the"
$ws.Range("C43").Value = "Heals the target, mana or health"

# Restore the view: scrolled so row 19 is at the top, with C29 selected.
$excel.Goto($ws.Range("A19"), $true) | Out-Null
$ws.Range("C29").Select() | Out-Null
